$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 16: et_var_name setting
$ws.Range("A16").Value = "et_var_name"
$ws.Range("B16").Value = "ET"
$ws.Range("C16").Value = "Choose between which variable to use for ET. ET is based on LE. ET_CORR is based energy-balance corrected (as done in OneFlux) LE_CORR variable"

# Copy formatting: A16/B16 like row 15 (A/B), C16 like row 6's C (wrap + vertical top, s=2)
$ws.Range("A15:B15").Copy() | Out-Null
$ws.Range("A16:B16").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$ws.Range("C6").Copy() | Out-Null
$ws.Range("C16").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$ws.Rows.Item(16).RowHeight = 72

# Row 17: blank row with same formatting as row 16 (columns A and B only)
$ws.Range("A16:B16").Copy() | Out-Null
$ws.Range("A17:B17").PasteSpecial(-4122) | Out-Null
$ws.Rows.Item(17).RowHeight = 23

$excel.CutCopyMode = 0

# Add data validation on B16 for ET / ET_CORR
$ws.Range("B16").Validation.Delete()
$ws.Range("B16").Validation.Add(3, 1, 1, """ET, ET_CORR""")
$ws.Range("B16").Validation.IgnoreBlank = $true
$ws.Range("B16").Validation.InCellDropdown = $true
$ws.Range("B16").Validation.ShowInput = $true
$ws.Range("B16").Validation.ShowError = $true

# Update view: scroll & selection (topLeftCell scrolling isn't round-tripped by
# this headless engine, but the active-cell selection is).
$ws.Application.ActiveWindow.ScrollRow = 6
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("A18").Select()
